# Applies the "Add big recip table. Use 1/x for lerp deltas." commit to
# P3DBenchmark/Performance.xlsx

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # "Tex, Flags<0>"
$ws2 = $wb.Worksheets.Item(2)   # "Flat, Flags<0>"
$ws3 = $wb.Worksheets.Item(3)   # "Tex, Flags<Subdivide = 16>"
$ws4 = $wb.Worksheets.Item(4)   # "Tex, Flags<Perspective>"
$ws5 = $wb.Worksheets.Item(5)   # "V1"

# ---------------------------------------------------------------------------
# Sheet1: "Tex, Flags<0>" - big recip summary table (rows 39-47) first, so
# the new shared strings "CRT Div".."RecipTable" land at indices 44-49
# BEFORE the "LerpDeltaRecip"/"DivRecip Lerp" labels (indices 50-51) used by
# the new per-sheet benchmark rows below.
# ---------------------------------------------------------------------------
$ws1.Range("B39").Value = "CRT Div"
$ws1.Range("C39").NumberFormat = "#,##0"
$ws1.Range("C39").Value = 124254

$ws1.Range("B40").Value = "FP Div"
$ws1.Range("C40").NumberFormat = "#,##0"
$ws1.Range("C40").Value = 124254

$ws1.Range("B41").Value = "FP Recip *"
$ws1.Range("C41").NumberFormat = "#,##0"
$ws1.Range("C41").Value = 1691188

$ws1.Range("C42").NumberFormat = "#,##0"
$ws1.Range("C43").NumberFormat = "#,##0"

$ws1.Range("B44").Value = "CRT recip"
$ws1.Range("C44").NumberFormat = "#,##0"
$ws1.Range("C44").Value = 160926

$ws1.Range("B45").Value = "FP Recip"
$ws1.Range("C45").NumberFormat = "#,##0"
$ws1.Range("C45").Value = 77869

$ws1.Range("B46").Value = "RecipTable"
$ws1.Range("C46").NumberFormat = "#,##0"
$ws1.Range("C46").Value = 160926

$ws1.Range("C47").NumberFormat = "#,##0"

# Extra little data points (columns K/L) added alongside the table.
$ws1.Range("K29").Value = 8.7279999999999998
$ws1.Range("K30").Value = 8.5890000000000004
$ws1.Range("K32").Value = 8.6069999999999993
$ws1.Range("L32").Value = 8.61

# New benchmark rows 20/21 ("LerpDeltaRecip" / "DivRecip Lerp") on sheet1.
$ws1.Range("A20").Value = "LerpDeltaRecip"
$ws1.Range("B20").Value = 92114
$ws1.Range("C20").Value = 97560
$ws1.Range("D20").Value = 116
$ws1.Range("H20").Value = 15320

$ws1.Range("A21").Value = "DivRecip Lerp"
$ws1.Range("B21").Value = 106598
$ws1.Range("C21").Value = 98648
$ws1.Range("D21").Value = 116
$ws1.Range("H21").Value = 15212

# H19 data point revised.
$ws1.Range("H19").Value = 15756

# Rebuild the E15:E21 shared-formula block (previously E15:E19 were each
# their own non-shared formula) and extend the F3:F19 block to F3:F21.
$ws1.Range("E15:E21").Formula = "=(D15/D`$2)-1"
$ws1.Range("F3:F21").Formula = "=(D3/D2)-1"

$ws1.Range("C22").Select()

# ---------------------------------------------------------------------------
# Sheet2: "Flat, Flags<0>" - append "LerpDeltaRecip"/"DivRecip Lerp" rows.
# ---------------------------------------------------------------------------
$ws2.Range("H13").Value = 15756

$ws2.Range("A14").Value = "LerpDeltaRecip"
$ws2.Range("B14").Value = 255297
$ws2.Range("C14").Value = 265957
$ws2.Range("D14").Value = 556
$ws2.Range("H14").Value = 15320

$ws2.Range("A15").Value = "DivRecip Lerp"
$ws2.Range("B15").Value = 253549
$ws2.Range("C15").Value = 276548
$ws2.Range("D15").Value = 549
$ws2.Range("H15").Value = 15212

$ws2.Range("E3:E15").Formula = "=(D3/D`$2)-1"
$ws2.Range("F3:F15").Formula = "=(D3/D2)-1"

$ws2.Range("A15").Select()

# ---------------------------------------------------------------------------
# Sheet3: "Tex, Flags<Subdivide = 16>" - append rows 12/13, plus some blank
# styled cells further down (I30:I35).
# ---------------------------------------------------------------------------
$ws3.Range("A12").Value = "LerpDeltaRecip"
$ws3.Range("B12").Value = 65316
$ws3.Range("C12").Value = 60328
$ws3.Range("D12").Value = 61
$ws3.Range("H12").Value = 16920

$ws3.Range("A13").Value = "DivRecip Lerp"
$ws3.Range("B13").Value = 43176
$ws3.Range("C13").Value = 41909
$ws3.Range("D13").Value = 58
$ws3.Range("H13").Value = 17012

$ws3.Range("E3:E13").Formula = "=(D3/D`$2)-1"
$ws3.Range("F4:F13").Formula = "=(D4/D3)-1"

$ws3.Range("I30").Style = $ws3.Range("E3").Style
$ws3.Range("I31").Style = $ws3.Range("E3").Style
$ws3.Range("I32").Style = $ws3.Range("E3").Style
$ws3.Range("I33").Style = $ws3.Range("E3").Style
$ws3.Range("I34").Style = $ws3.Range("E3").Style
$ws3.Range("I35").Style = $ws3.Range("E3").Style

$ws3.Range("A13").Select()

# ---------------------------------------------------------------------------
# Sheet4: "Tex, Flags<Perspective>" - append rows 10/11.
# ---------------------------------------------------------------------------
$ws4.Range("A10").Value = "LerpDeltaRecip"
$ws4.Range("B10").Value = 40561
$ws4.Range("C10").Value = 32744
$ws4.Range("D10").Value = 31
$ws4.Range("H10").Value = 19668

$ws4.Range("A11").Value = "DivRecip Lerp"
$ws4.Range("B11").Value = 17938
$ws4.Range("C11").Value = 17419
$ws4.Range("D11").Value = 21
$ws4.Range("H11").Value = 19568

$ws4.Range("E3:E11").Formula = "=(D3/D`$2)-1"
$ws4.Range("F7:F11").Formula = "=(D7/D6)-1"

$ws4.Range("C11").Select()

# ---------------------------------------------------------------------------
# Sheet5: "V1" - literal data point refresh (no structural/formula changes,
# cached values follow from the existing formulas once the recalc runs).
# ---------------------------------------------------------------------------
$ws5.Range("D2").Value = 116
$ws5.Range("H2").Value = 15212
$ws5.Range("D6").Value = 549
$ws5.Range("D9").Value = 58
$ws5.Range("H9").Value = 17012

$ws5.Range("E9").Select()

Write-Output "edit applied"
